# Reorder the comma-separated "Net Name" lists in the BoM and DNF sheets.
# (Same set of net names per row, just written back in a different order,
# matching the upstream KiBot regeneration referenced by the commit message.)

$wb = $excel.ActiveWorkbook

$wsBoM = $wb.Worksheets.Item("BoM")
$wsDNF = $wb.Worksheets.Item("DNF")

# --- BoM sheet (Net Name column = X) ---
$wsBoM.Range("X10").Value = "Earth,Net-(U1-UCAP)"
$wsBoM.Range("X13").Value = "Net-(D2-A),/RXLED"
$wsBoM.Range("X14").Value = ",GND,+5V"
$wsBoM.Range("X15").Value = "/SCK2,/MOSI2,/RESET2,/MISO2,GND,+5V"
$wsBoM.Range("X16").Value = "Net-(J3-Pin_5),Net-(J3-Pin_3),Net-(J3-Pin_1),Net-(J3-Pin_2),Net-(J3-Pin_4)"
$wsBoM.Range("X17").Value = "Net-(J6-Pin_6),Net-(J6-Pin_5),Net-(J6-Pin_4),Net-(J6-Pin_1),Net-(J6-Pin_2),Net-(J6-Pin_3)"
$wsBoM.Range("X21").Value = "/RESET2,Net-(J4-Pin_1),Net-(U1-XTAL1),Net-(J4-Pin_2),/MISO2,/MOSI2,Net-(J3-Pin_3),Net-(U1-UCAP),/DTR,VBUS,Net-(J4-Pin_3),Net-(J6-Pin_3),Net-(J3-Pin_4),Net-(U1-PC0{slash}XTAL2),+5V,/SCK2,Net-(J3-Pin_5),Net-(J6-Pin_4),Net-(U1-D+),Net-(J4-Pin_4),Net-(J3-Pin_2),Net-(U1-D-),GND,unconnected-(U1-PB0-Pad14),Net-(J6-Pin_2),Net-(J6-Pin_6),Net-(J6-Pin_5),Net-(J3-Pin_1),/RXLED,/TXLED,Earth"

# --- DNF sheet (Net Name column = X) ---
$wsDNF.Range("X9").Value = "GND,Net-(U1-XTAL1)"
$wsDNF.Range("X12").Value = "Net-(J4-Pin_4),Net-(J4-Pin_2),Net-(J4-Pin_1),Net-(J4-Pin_3)"
$wsDNF.Range("X13").Value = "Net-(J2-VBUS),Earth,Net-(J2-Shield),Net-(J2-D+),Net-(J2-D-)"
$wsDNF.Range("X15").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
$wsDNF.Range("X16").Value = "Net-(J2-D+),Net-(J2-Shield)"
# X17 shares the same original text as X15 ("Net-(U1-XTAL1),Net-(U1-PC0{slash}XTAL2)")
# and must end up with the same new text too.
$wsDNF.Range("X17").Value = "Net-(U1-PC0{slash}XTAL2),Net-(U1-XTAL1)"
